# CCO_eCoaching_Log_Admin_Runbook.docx
#
# "Get file from TFS" step references the changeset number of the
# ecl_admin_publish.zip build, shown in large bold text right after the
# literal word "Changeset ". Bump it from 51959 to 52007 (the
# "Updated zip changeset#" commit).
#
# "51959" only occurs once in the whole document, so a plain
# Find/Replace on just the digits is enough to hit the right spot
# without disturbing the "Changeset " label (which uses different,
# smaller, non-bold formatting) or any of the surrounding runs.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "51959",  # FindText
    $true,    # MatchCase
    $false,   # MatchWholeWord
    $false,   # MatchWildcards
    $false,   # MatchSoundsLike
    $false,   # MatchAllWordForms
    $true,    # Forward
    1,        # Wrap (wdFindContinue)
    $false,   # Format
    "52007",  # ReplaceWith
    2         # Replace (wdReplaceAll)
)
